$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 63, pushing existing rows 63-79 down to 64-80
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new data entry
$ws.Cells.Item(63, 1).Value = 1
$ws.Cells.Item(63, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(63, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(63, 4).Value = 44784
$ws.Cells.Item(63, 4).NumberFormat = $ws.Cells.Item(64, 4).NumberFormat
$ws.Cells.Item(63, 5).Value = 15
$ws.Cells.Item(63, 6).Value = 100114001
$ws.Cells.Item(63, 7).Value = "Papa"
$ws.Cells.Item(63, 8).Value = "Asterix"
$ws.Cells.Item(63, 9).Value = "1a (guarda)"
$ws.Cells.Item(63, 10).Value = 1000
$ws.Cells.Item(63, 11).Value = 9000
$ws.Cells.Item(63, 12).Value = 10000
$ws.Cells.Item(63, 13).Value = 9500
$ws.Cells.Item(63, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(63, 16).Value = 380
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
